$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 6 with the BTMR_LEVEL_2 role data
$ws.Range("A6").Value = "BTMR_LEVEL_2"
$ws.Range("B6").Value = "Pentadbir Helpdesk Sistem Second Level (BTMR)"
$ws.Range("C6").Value = "System Helpdesk Administrator Second Level (BTMR)"
$ws.Range("D6").Value = "Mengurus konfigurasi sistem helpdesk (BTMR) termasuk modul, kategori, pengguna, dan peranan."

# Update selection to the new last row, mirroring the original file's pattern
$ws.Range("A6").Select()

$wb.Save()
